$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

$ws.Range("B12").Value = 75
$ws.Range("C12").Value = -14.4
$ws.Range("E12").Value = "60.6/140"
